# Update res_bus/vm_pu.xlsx "Sheet1" data for the 380 kV case (Case_2_244):
# - B2:B25 (slack bus vm_pu setpoint) drops from 1.05 to 1.02 p.u.
# - C2:F25 and I2:N25 (other bus voltages) are refreshed to the recomputed
#   load-flow results that follow from the new slack voltage.
# - Columns G (always 1) and H (always empty) are left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$arr = New-Object 'object[,]' 24,13
$arr[0,0] = 1.02
$arr[0,1] = 1.035551122635881
$arr[0,2] = 1.042561351605305
$arr[0,3] = 0.992614727750844
$arr[0,4] = 1.051532467930445
$arr[0,5] = 1
$arr[0,6] = $null
$arr[0,7] = 1.038245119272223
$arr[0,8] = 1.04066428011297
$arr[0,9] = 1.045337691070472
$arr[0,10] = 0.9955398523335997
$arr[0,11] = 1.054283708423865
$arr[0,12] = 1.017477945689646
$arr[1,0] = 1.02
$arr[1,1] = 1.036687812631788
$arr[1,2] = 1.043423019422623
$arr[1,3] = 0.9936372048519299
$arr[1,4] = 1.052537585810179
$arr[1,5] = 1
$arr[1,6] = $null
$arr[1,7] = 1.038504504137332
$arr[1,8] = 1.04144362110775
$arr[1,9] = 1.046010082430127
$arr[1,10] = 0.9963617723202687
$arr[1,11] = 1.055100996538723
$arr[1,12] = 1.017738224862476
$arr[2,0] = 1.02
$arr[2,1] = 1.03742294093165
$arr[2,2] = 1.043979995198189
$arr[2,3] = 0.9942998659930998
$arr[2,4] = 1.053187665690316
$arr[2,5] = 1
$arr[2,6] = $null
$arr[2,7] = 1.03867043755345
$arr[2,8] = 1.041946992407752
$arr[2,9] = 1.046443945411629
$arr[2,10] = 0.9968940712668347
$arr[2,11] = 1.055628927746717
$arr[2,12] = 1.01790626770715
$arr[3,0] = 1.02
$arr[3,1] = 1.037731897598503
$arr[3,2] = 1.044214008775504
$arr[3,3] = 0.994578699834602
$arr[3,4] = 1.053460887968347
$arr[3,5] = 1
$arr[3,6] = $null
$arr[3,7] = 1.038739739638141
$arr[3,8] = 1.042158391217426
$arr[3,9] = 1.046626049694986
$arr[3,10] = 0.9971179600053012
$arr[3,11] = 1.05585065262693
$arr[3,12] = 1.017976823205278
$arr[4,0] = 1.02
$arr[4,1] = 1.037783767464148
$arr[4,2] = 1.044253292532377
$arr[4,3] = 0.994625531979634
$arr[4,4] = 1.053506759014896
$arr[4,5] = 1
$arr[4,6] = $null
$arr[4,7] = 1.038751349002843
$arr[4,8] = 1.042193873208168
$arr[4,9] = 1.046656608703763
$arr[4,10] = 0.9971555583673455
$arr[4,11] = 1.055887868466645
$arr[4,12] = 1.017988664526057
$arr[5,0] = 1.02
$arr[5,1] = 1.037427069582931
$arr[5,2] = 1.043983122645371
$arr[5,3] = 0.994303590798249
$arr[5,4] = 1.053191316779737
$arr[5,5] = 1
$arr[5,6] = $null
$arr[5,7] = 1.038671365364402
$arr[5,8] = 1.041949817986966
$arr[5,9] = 1.046446379844033
$arr[5,10] = 0.9968970624462089
$arr[5,11] = 1.055631891299448
$arr[5,12] = 1.017907210825556
$arr[6,0] = 1.02
$arr[6,1] = 1.035935353001038
$arr[6,2] = 1.042852676586276
$arr[6,3] = 0.9929600610674297
$arr[6,4] = 1.051872214163184
$arr[6,5] = 1
$arr[6,6] = $null
$arr[6,7] = 1.038333174217528
$arr[6,8] = 1.040927851654514
$arr[6,9] = 1.045565181420701
$arr[6,10] = 0.9958175282591056
$arr[6,11] = 1.054560103065317
$arr[6,12] = 1.017565986003915
$arr[7,0] = 1.02
$arr[7,1] = 1.033303747291741
$arr[7,2] = 1.040856234387264
$arr[7,3] = 0.9906006454969559
$arr[7,4] = 1.04954549069018
$arr[7,5] = 1
$arr[7,6] = $null
$arr[7,7] = 1.037722642065231
$arr[7,8] = 1.039119992681582
$arr[7,9] = 1.044003051460288
$arr[7,10] = 0.9939188001724441
$arr[7,11] = 1.052664505380456
$arr[7,12] = 1.016961825762714
$arr[8,0] = 1.02
$arr[8,1] = 1.031547241863367
$arr[8,2] = 1.039522272311946
$arr[8,3] = 0.989033133672735
$arr[8,4] = 1.047992781108499
$arr[8,5] = 1
$arr[8,6] = $null
$arr[8,7] = 1.037305806336327
$arr[8,8] = 1.037909999658138
$arr[8,9] = 1.042955333246334
$arr[8,10] = 0.9926553831429383
$arr[8,11] = 1.051396071088834
$arr[8,12] = 1.016557108789284
$arr[9,0] = 1.02
$arr[9,1] = 1.030786138898844
$arr[9,2] = 1.038943937482413
$arr[9,3] = 0.988355674866747
$arr[9,4] = 1.047320065611886
$arr[9,5] = 1
$arr[9,6] = $null
$arr[9,7] = 1.037122983217668
$arr[9,8] = 1.037384923874372
$arr[9,9] = 1.042500161572177
$arr[9,10] = 0.9921088820399291
$arr[9,11] = 1.050845704926958
$arr[9,12] = 1.01638139875721
$arr[10,0] = 1.02
$arr[10,1] = 1.030503350875613
$arr[10,2] = 1.038729009566249
$arr[10,3] = 0.9881042295826724
$arr[10,4] = 1.047070131173244
$arr[10,5] = 1
$arr[10,6] = $null
$arr[10,7] = 1.037054724278529
$arr[10,8] = 1.03718971515754
$arr[10,9] = 1.042330864285878
$arr[10,10] = 0.9919059725120875
$arr[10,11] = 1.050641104764227
$arr[10,12] = 1.016316062154437
$arr[11,0] = 1.02
$arr[11,1] = 1.03056401353631
$arr[11,2] = 1.038775117254768
$arr[11,3] = 0.9881581567098651
$arr[11,4] = 1.047123745590235
$arr[11,5] = 1
$arr[11,6] = $null
$arr[11,7] = 1.03706938190683
$arr[11,8] = 1.037231595894875
$arr[11,9] = 1.042367189375119
$arr[11,10] = 0.9919494934313052
$arr[11,11] = 1.050684999880071
$arr[11,12] = 1.016330080252143
$arr[12,0] = 1.02
$arr[12,1] = 1.030762765204491
$arr[12,2] = 1.038926173691614
$arr[12,3] = 0.9883348863814464
$arr[12,4] = 1.047299407138666
$arr[12,5] = 1
$arr[12,6] = $null
$arr[12,7] = 1.037117348059287
$arr[12,8] = 1.03736879137824
$arr[12,9] = 1.042486172024155
$arr[12,10] = 0.9920921077337197
$arr[12,11] = 1.050828796077929
$arr[12,12] = 1.016375999444033
$arr[13,0] = 1.02
$arr[13,1] = 1.030885211949926
$arr[13,2] = 1.039019230144667
$arr[13,3] = 0.9884438009545853
$arr[13,4] = 1.047407630320295
$arr[13,5] = 1
$arr[13,6] = $null
$arr[13,7] = 1.037146855166776
$arr[13,8] = 1.037453299202511
$arr[13,9] = 1.042559451164824
$arr[13,10] = 0.9921799884222134
$arr[13,11] = 1.050917371162845
$arr[13,12] = 1.016404282481483
$arr[14,0] = 1.02
$arr[14,1] = 1.031597742565306
$arr[14,2] = 1.039560639279465
$arr[14,3] = 0.9890781214508737
$arr[14,4] = 1.048037418924671
$arr[14,5] = 1
$arr[14,6] = $null
$arr[14,7] = 1.037317890590031
$arr[14,8] = 1.03794482310507
$arr[14,9] = 1.042985509788367
$arr[14,10] = 0.9926916645766087
$arr[14,11] = 1.051432573314872
$arr[14,12] = 1.016568760271074
$arr[15,0] = 1.02
$arr[15,1] = 1.032044552575862
$arr[15,2] = 1.039900057783488
$arr[15,3] = 0.989476357848556
$arr[15,4] = 1.04843236601394
$arr[15,5] = 1
$arr[15,6] = $null
$arr[15,7] = 1.037424552539976
$arr[15,8] = 1.038252837074339
$arr[15,9] = 1.043252362432274
$arr[15,10] = 0.9930127773699352
$arr[15,11] = 1.051755444245967
$arr[15,12] = 1.016671808214679
$arr[16,0] = 1.02
$arr[16,1] = 1.032305118800656
$arr[16,2] = 1.040097965367111
$arr[16,3] = 0.9897087662937556
$arr[16,4] = 1.04866269514969
$arr[16,5] = 1
$arr[16,6] = $null
$arr[16,7] = 1.037486541802807
$arr[16,8] = 1.038432386392792
$arr[16,9] = 1.043407868101914
$arr[16,10] = 0.9932001317071769
$arr[16,11] = 1.051943660816119
$arr[16,12] = 1.016731869494798
$arr[17,0] = 1.02
$arr[17,1] = 1.032393956636142
$arr[17,2] = 1.040165434963602
$arr[17,3] = 0.9897880325774034
$arr[17,4] = 1.048741225204362
$arr[17,5] = 1
$arr[17,6] = $null
$arr[17,7] = 1.037507640406193
$arr[17,8] = 1.038493589433436
$arr[17,9] = 1.04346086694305
$arr[17,10] = 0.9932640239640975
$arr[17,11] = 1.052007819401412
$arr[17,12] = 1.016752341240042
$arr[18,0] = 1.02
$arr[18,1] = 1.031996619325156
$arr[18,2] = 1.039863648580116
$arr[18,3] = 0.9894336180360679
$arr[18,4] = 1.048389995741241
$arr[18,5] = 1
$arr[18,6] = $null
$arr[18,7] = 1.037413131980334
$arr[18,8] = 1.038219801475751
$arr[18,9] = 1.043223746666558
$arr[18,10] = 0.9929783193494215
$arr[18,11] = 1.051720814485435
$arr[18,12] = 1.016660756785345
$arr[19,0] = 1.02
$arr[19,1] = 1.030704240027618
$arr[19,2] = 1.03888169433715
$arr[19,3] = 0.9882828385668249
$arr[19,4] = 1.047247680795166
$arr[19,5] = 1
$arr[19,6] = $null
$arr[19,7] = 1.037103232894035
$arr[19,8] = 1.037328395487672
$arr[19,9] = 1.042451140849407
$arr[19,10] = 0.9920501090198102
$arr[19,11] = 1.050786456369623
$arr[19,12] = 1.016362479322721
$arr[20,0] = 1.02
$arr[20,1] = 1.029891203311606
$arr[20,2] = 1.038263672154656
$arr[20,3] = 0.9875604150241495
$arr[20,4] = 1.046529126481516
$arr[20,5] = 1
$arr[20,6] = $null
$arr[20,7] = 1.036906360268727
$arr[20,8] = 1.036766936514173
$arr[20,9] = 1.041964063589778
$arr[20,10] = 0.9914670000341481
$arr[20,11] = 1.050198006435066
$arr[20,12] = 1.016174534906733
$arr[21,0] = 1.02
$arr[21,1] = 1.030322254249152
$arr[21,2] = 1.038591357115625
$arr[21,3] = 0.9879432794643023
$arr[21,4] = 1.046910077716442
$arr[21,5] = 1
$arr[21,6] = $null
$arr[21,7] = 1.037010918342434
$arr[21,8] = 1.037064671278398
$arr[21,9] = 1.042222396733789
$arr[21,10] = 0.991776070289318
$arr[21,11] = 1.050510048133715
$arr[21,12] = 1.016274206313981
$arr[22,0] = 1.02
$arr[22,1] = 1.032018278453408
$arr[22,2] = 1.039880100546374
$arr[22,3] = 0.9894529299347244
$arr[22,4] = 1.048409141157152
$arr[22,5] = 1
$arr[22,6] = $null
$arr[22,7] = 1.037418293134257
$arr[22,8] = 1.038234729180388
$arr[22,9] = 1.043236677347002
$arr[22,10] = 0.9929938892766442
$arr[22,11] = 1.051736462518831
$arr[22,12] = 1.016665750589018
$arr[23,0] = 1.02
$arr[23,1] = 1.033984445157039
$arr[23,2] = 1.04137289083614
$arr[23,3] = 0.9912096547607049
$arr[23,4] = 1.050147278269126
$arr[23,5] = 1
$arr[23,6] = $null
$arr[23,7] = 1.037882209021786
$arr[23,8] = 1.039588203100504
$arr[23,9] = 1.044408008491422
$arr[23,10] = 0.9944092447426414
$arr[23,11] = 1.053155390427034
$arr[23,12] = 1.017118357610822
$ws.Range("B2:N25").Value = $arr
Write-Host "applied vm_pu updates for the 380 kV case"
